# Applies the two content edits described by the diff:
#   1. A new run containing a single back-tick ("`") is inserted immediately
#      before the existing "Assignment" run in the first paragraph. The new
#      run shares the same (Bold, sz 40) formatting but is a distinct
#      <w:r> element (it carries no w:rsidRPr, unlike the original run).
#   2. The three runs that host the inline pictures (rId7, rId8, rId9) gain
#      a <w:noProof/> flag in their run properties, right after <w:bCs/>.
#
# (The third hunk in the source diff only reorders the xmlns:* attributes
#  on the pre-existing a14:useLocalDpi extension element -- a byte-level
#  artifact of a real Word round-trip save. It carries no semantic meaning
#  -- the attribute set/values are identical either way -- and it is not
#  reachable through any Word object-model property, so it is intentionally
#  left untouched here.)

$d = $word.ActiveDocument

# --- Change 1: insert the back-tick run before "Assignment" -----------------
$p1 = $d.Paragraphs(1)
$r = $p1.Range

$xmlFrag = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
  '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
    '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
      '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
          '<w:body>' +
            '<w:p w:rsidR="00AF59FA" w:rsidRPr="00DD2DEB" w:rsidRDefault="00DF53A6" w:rsidP="00DF53A6">' +
              '<w:pPr><w:jc w:val="center"/><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr></w:pPr>' +
              '<w:r><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>`</w:t></w:r>' +
              '<w:r w:rsidRPr="00DD2DEB"><w:rPr><w:b/><w:bCs/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t>Assignment</w:t></w:r>' +
            '</w:p>' +
          '</w:body>' +
        '</w:document>' +
      '</pkg:xmlData>' +
    '</pkg:part>' +
  '</pkg:package>'

$r.InsertXML($xmlFrag)

# --- Change 2: flag the three picture runs as NoProofing (<w:noProof/>) ----
for ($i = 1; $i -le $d.InlineShapes.Count; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = $true
}
